# "Corrected Lookup xl file"
#
# The Lookup sheet (Sheet1) lists account owners in column B with a
# lookup "Type" flag in column C. Row 2 (BankOfThePhilippineIslands)
# already carried a numeric marker of 1 in column A (formatted as a
# three-digit "000" code). This fix back-fills the same marker for the
# remaining rows (3-11) so every entry is consistently flagged.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Back-fill column A (rows 3 through 11) with the same "1" marker /
# number format ("000") already used on A2, matching the existing
# style so no new cell style is introduced.
for ($r = 3; $r -le 11; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value = 1
    $cell.NumberFormat = "000"
}

# Move the sheet's active selection from B14 to D14.
$ws.Range("D14").Select()

# Restore the workbook window to its normal (non-maximized) size and
# position, as last saved by the author.
$window = $excel.Windows.Item(1)
$window.Left = 705
$window.Top = 690
$window.Width = 28800
$window.Height = 11295
